$p = $ppt.ActivePresentation

# The deck's single theme (ppt/theme/theme2.xml, wired to the slide master)
# currently holds the "Integral" / "Red Violet" color scheme. The edit swaps
# the "Office Theme" colors (previously sitting unused in ppt/theme/theme1.xml,
# only referenced by the notes master) into the master/slide theme, and the
# "Red Violet" colors into the other theme slot. Only the color values can be
# pushed through the COM surface, so we overwrite the 12 scheme colors here
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - the exact element order
# used inside <a:clrScheme>) with the values that used to live in theme1.xml.

function New-ComRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    New-ComRgb 0x00 0x00 0x00,   # dk1
    New-ComRgb 0xFF 0xFF 0xFF,   # lt1
    New-ComRgb 0x44 0x54 0x6A,   # dk2
    New-ComRgb 0xE7 0xE6 0xE6,   # lt2
    New-ComRgb 0x5B 0x9B 0xD5,   # accent1
    New-ComRgb 0xED 0x7D 0x31,   # accent2
    New-ComRgb 0xA5 0xA5 0xA5,   # accent3
    New-ComRgb 0xFF 0xC0 0x00,   # accent4
    New-ComRgb 0x44 0x72 0xC4,   # accent5
    New-ComRgb 0x70 0xAD 0x47,   # accent6
    New-ComRgb 0x05 0x63 0xC1,   # hlink
    New-ComRgb 0x95 0x4F 0x72    # folHlink
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
